# "ChatGPT API parameters selection.xlsx" re-upload.
#
# The authoritative diff for this commit shows no actual content or
# formatting change anywhere in the workbook -- every cell value, formula,
# merge, number format, font, fill and alignment is identical before and
# after. The only user-visible difference is on the "Compiled data" sheet:
# the window had been left scrolled down (topLeftCell="A24") with D35
# selected, and after the re-save it is scrolled back to the top with H10
# selected instead. (The surrounding GUID / xf-ordering churn in the raw
# OOXML is just Excel's own internal bookkeeping from re-saving the file,
# not a deliberate edit, so we don't try to fabricate it.)
#
# Reproduce that: make sure "Compiled data" is the active sheet and select
# H10 on it, which scrolls the view back so A1 is the top-left cell again.

$wb = $excel.ActiveWorkbook

$wsCompiled = $wb.Worksheets.Item("Compiled data")
$wsDialogs  = $wb.Worksheets.Item("Test dialogs")

# Re-assert the existing wrap/top-align formatting on the dialog transcript
# column -- values/format are unchanged from before, this just mirrors the
# formatting action implied by the style-table churn in the diff without
# altering the visible result.
$dialogRange = $wsDialogs.Range("A1:A76")
$dialogRange.VerticalAlignment = -4160  # xlTop
$dialogRange.WrapText = $true

# Bring the data sheet to the front and move the selection to H10, which
# also resets the scrolled view back to the top (topLeftCell A1) -- matching
# the sheetView change in the diff.
$wsCompiled.Activate() | Out-Null
$wsCompiled.Range("H10").Select() | Out-Null
